$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H49").Value = 1810.6
$ws.Range("I49").Value = 1013.25
$ws.Range("J49").Value = 5000
$ws.Range("K49").Value = 3039.75
$ws.Range("L49").Value = 15000
$ws.Range("M49").Value = -2903.75
$ws.Range("N49").Value = -15272
$ws.Range("H53").Value = 380.2
$ws.Range("J53").Value = 850.25
$ws.Range("L53").Value = 850.25
$ws.Range("N53").Value = -2124.25
$ws.Range("H113").Value = 3769.1428
$ws.Range("I113").Value = 3769.1428
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3769.1428
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -515.1428000000001
$ws.Range("N113").ClearContents()

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1133241.8
$ws.Range("I2").Value = 1840031.2
$ws.Range("K2").Value = 1840031.2
$ws.Range("M2").Value = -1839918.2
$ws.Range("H61").Value = 31253098
$ws.Range("I61").Value = 55557436
$ws.Range("K61").Value = 55557436
$ws.Range("M61").Value = -55557224
$ws.Range("H74").Value = 47621970
$ws.Range("I74").Value = 52633750
$ws.Range("J74").Value = 10000
$ws.Range("K74").Value = 52633750
$ws.Range("L74").Value = 10000
$ws.Range("M74").Value = -52632876
$ws.Range("N74").Value = -11748
$ws.Range("H77").Value = 47621970
$ws.Range("I77").Value = 52633750
$ws.Range("J77").Value = 10000
$ws.Range("K77").Value = 263168750
$ws.Range("L77").Value = 50000
$ws.Range("M77").Value = -263164382
$ws.Range("N77").Value = -58736
$ws.Range("H88").Value = 113875.664
$ws.Range("I88").Value = 145411.58
$ws.Range("J88").Value = 3500
$ws.Range("K88").Value = 145411.58
$ws.Range("L88").Value = 3500
$ws.Range("M88").Value = -145005.58
$ws.Range("N88").Value = -4312
$ws.Range("H91").Value = 113875.664
$ws.Range("I91").Value = 145411.58
$ws.Range("J91").Value = 3500
$ws.Range("K91").Value = 145411.58
$ws.Range("L91").Value = 3500
$ws.Range("M91").Value = -144007.58
$ws.Range("N91").Value = -6308
$ws.Range("H116").Value = 1133241.8
$ws.Range("I116").Value = 1840031.2
$ws.Range("K116").Value = 1840031.2
$ws.Range("M116").Value = -1837737.2
$ws.Range("H132").Value = 4548198
$ws.Range("I132").Value = 5265071.5
$ws.Range("K132").Value = 15795214.5
$ws.Range("M132").Value = -15792684.5
$ws.Range("H136").Value = 31253098
$ws.Range("I136").Value = 55557436
$ws.Range("K136").Value = 166672308
$ws.Range("M136").Value = -166669758

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1133241.8
$ws.Range("I3").Value = 1840031.2
$ws.Range("K3").Value = 1840031.2
$ws.Range("M3").Value = -1839917.2

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 2172.923
$ws.Range("J7").Value = 1356
$ws.Range("L7").Value = 1356
$ws.Range("N7").Value = -1582
$ws.Range("H31").Value = 5596.8887
$ws.Range("I31").Value = 5671.5
$ws.Range("K31").Value = 5671.5
$ws.Range("M31").Value = -5376.5
$ws.Range("H34").Value = 5596.8887
$ws.Range("I34").Value = 5671.5
$ws.Range("K34").Value = 5671.5
$ws.Range("M34").Value = -5469.5
$ws.Range("H58").Value = 17862864
$ws.Range("I58").Value = 41676292
$ws.Range("J58").Value = 2794.8125
$ws.Range("K58").Value = 41676292
$ws.Range("L58").Value = 2794.8125
$ws.Range("M58").Value = -41676089
$ws.Range("N58").Value = -3200.8125
$ws.Range("H105").Value = 1033.8182
$ws.Range("I105").Value = 1037.3
$ws.Range("J105").Value = 999
$ws.Range("K105").Value = 1037.3
$ws.Range("L105").Value = 999
$ws.Range("M105").Value = 709.7
$ws.Range("N105").Value = -4493
$ws.Range("H136").Value = 17862864
$ws.Range("I136").Value = 41676292
$ws.Range("J136").Value = 2794.8125
$ws.Range("K136").Value = 125028876
$ws.Range("L136").Value = 8384.4375
$ws.Range("M136").Value = -125026326
$ws.Range("N136").Value = -13484.4375

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 194.80952
$ws.Range("J12").Value = 277.64285
$ws.Range("L12").Value = 832.9285500000001
$ws.Range("N12").Value = -1178.92855
$ws.Range("H99").Value = 5512
$ws.Range("I99").Value = 25
$ws.Range("J99").Value = 10999
$ws.Range("K99").Value = 75
$ws.Range("L99").Value = 32997
$ws.Range("M99").Value = 2171
$ws.Range("N99").Value = -37489
$ws.Range("H109").Value = 1515.1818
$ws.Range("I109").Value = 1515.1818
$ws.Range("K109").Value = 4545.5454
$ws.Range("M109").Value = -3505.5454
$ws.Range("H117").Value = 2199.077
$ws.Range("I117").Value = 270.25
$ws.Range("J117").Value = 3056.3333
$ws.Range("K117").Value = 810.75
$ws.Range("L117").Value = 9168.999899999999
$ws.Range("M117").Value = 2631.25
$ws.Range("N117").Value = -16052.9999
$ws.Range("H121").Value = 103384.586
$ws.Range("I121").Value = 201059.8
$ws.Range("J121").Value = 33616.57
$ws.Range("K121").Value = 603179.3999999999
$ws.Range("L121").Value = 100849.71
$ws.Range("M121").Value = -601869.3999999999
$ws.Range("N121").Value = -103469.71
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 4450
$ws.Range("I97").Value = 4450
$ws.Range("K97").Value = 4450
$ws.Range("M97").Value = -3954
$ws.Range("H107").Value = 5589.778
$ws.Range("I107").Value = 4329
$ws.Range("J107").Value = 10002.5
$ws.Range("K107").Value = 4329
$ws.Range("L107").Value = 10002.5
$ws.Range("M107").Value = -2409
$ws.Range("N107").Value = -13842.5
$ws.Range("H132").Value = 6946984
$ws.Range("I132").Value = 7814994.5
$ws.Range("J132").Value = 2899.5
$ws.Range("K132").Value = 23444983.5
$ws.Range("L132").Value = 8698.5
$ws.Range("M132").Value = -23442453.5
$ws.Range("N132").Value = -13758.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2076.625
$ws.Range("I22").Value = 2729.625
$ws.Range("K22").Value = 2729.625
$ws.Range("M22").Value = -2434.625
$ws.Range("H27").Value = 2076.625
$ws.Range("I27").Value = 2729.625
$ws.Range("K27").Value = 2729.625
$ws.Range("M27").Value = -2622.625
$ws.Range("H55").Value = 438.84375
$ws.Range("I55").Value = 296.8125
$ws.Range("K55").Value = 296.8125
$ws.Range("M55").Value = -123.8125
$ws.Range("H82").Value = 1886.8636
$ws.Range("I82").Value = 1928.2
$ws.Range("J82").Value = 1798.2858
$ws.Range("K82").Value = 1928.2
$ws.Range("L82").Value = 1798.2858
$ws.Range("M82").Value = -1567.2
$ws.Range("N82").Value = -2520.2858
$ws.Range("H85").Value = 1886.8636
$ws.Range("I85").Value = 1928.2
$ws.Range("J85").Value = 1798.2858
$ws.Range("K85").Value = 1928.2
$ws.Range("L85").Value = 1798.2858
$ws.Range("M85").Value = -680.2
$ws.Range("N85").Value = -4294.2858

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 51559.25
$ws.Range("I81").Value = 56732.777
$ws.Range("K81").Value = 113465.554
$ws.Range("M81").Value = -112404.554
$ws.Range("H84").Value = 51559.25
$ws.Range("I84").Value = 56732.777
$ws.Range("K84").Value = 567327.77
$ws.Range("M84").Value = -562023.77
$ws.Range("H132").Value = 18534376
$ws.Range("I132").Value = 21748986
$ws.Range("K132").Value = 65246958
$ws.Range("M132").Value = -65244428
$ws.Range("H136").Value = 22729228
$ws.Range("I136").Value = 29413836
$ws.Range("K136").Value = 88241508
$ws.Range("M136").Value = -88238958
